# Adding Sign up scenario in driver app
# Adds a new "url.driverapp.qa" / driver-app URL row to the URL sheet,
# mirroring the formatting of the existing rows and wiring up a new
# hyperlink relationship for the URL cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URL")

# --- Mark the (previously blank) filler column C for the existing rows ---
# (matches the CREDENTIALS sheet pattern of a trailing style-only column)
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C4").NumberFormat = "@"

# --- New row 5: url.driverapp.qa / https://driverapp.dev.mobile22.com ---
# Set A5 first so its shared-string entry is interned before B5's, matching
# the original authoring order (url.driverapp.qa, then the URL string).
$ws.Range("A5").Value = "url.driverapp.qa"
$ws.Range("B5").Value = "https://driverapp.dev.mobile22.com"

# Hyperlink the new URL cell first (this also nudges its style, which gets
# normalised back to match the other rows by the formatting copy below).
$ws.Hyperlinks.Add($ws.Range("B5"), "https://driverapp.dev.mobile22.com")

# Clone the formatting of row 4 (A4:B4) onto row 5, so the new cells pick up
# the exact same styles already used by the other data rows.
$ws.Range("A4:B4").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)

$ws.Range("A5").Value = "url.driverapp.qa"
$ws.Range("B5").Value = "https://driverapp.dev.mobile22.com"

$ws.Range("C5").NumberFormat = "@"

# --- New trailing row 6: just the style-only filler column C ---
$ws.Range("C6").NumberFormat = "@"

# Match the author's final cursor position.
[void]$ws.Range("D15").Select()
